$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("B10").Value = "[Leonardo R.-Retificação, -, -, -]"
$ws.Range("C10").Value = "-"
$ws.Range("D10").Value = "[-, Ludoff-Cont. Lóg. Prog. CLP, Leonardo R.-Mec. Manut. Equip. Ind., Anselmo-M. Motor Endot.]"
$ws.Range("E10").Value = "[Humberto-Comam. Pneumáticos, Valmir-Calderaria, Gisele-Ens. Dest. não Dest., Cleidson-Metrologia 2]"
$ws.Range("F10").Value = "Rogério-Elem"

# Row 11
$ws.Range("B11").Value = "[Leonardo R.-Retificação, -, -, -]"
$ws.Range("C11").Value = "-"
$ws.Range("D11").Value = "[Leonardo R.-Retificação, -, Ludoff-Camam. Hidráulicos, Anselmo-M. Motor Endot.]"
$ws.Range("E11").Value = "[Humberto-Comam. Pneumáticos, Valmir-Calderaria, Gisele-Ens. Dest. não Dest., Cleidson-Metrologia 2]"
$ws.Range("F11").Value = "Rogério-Elem"

# Row 12
$ws.Range("B12").Value = "[Ivan-Tec. Soldagem, Aderci-M. A. Comp; Cad / CAM, And. Edson-Usin. CNC, Ludoff-Cont. Lóg. Prog. CLP]"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "[Aderci-Fresagem, Anselmo-M. Motor Endot., Ludoff-Cont. Lóg. Prog. CLP, Joel Lima-Tec. Fundição]"
$ws.Range("E12").Value = "[Humberto-Comam. Pneumáticos, Valmir-Calderaria, Gisele-Ens. Dest. não Dest., Cleidson-Metrologia 2]"

# Row 14
$ws.Range("B14").Value = "[Ivan-Tec. Soldagem, Wellington-Trat. Térmicos, And. Edson-Usin. CNC, Aderci-Fresagem]"
$ws.Range("C14").Value = "[-, -, Leonardo R.-Mec. Manut. Equip. Ind., -]"
$ws.Range("D14").Value = "[Joel Lima-Tec. Fundição, Ludoff-Camam. Hidráulicos, And. Edson-Usin. CNC, Aderci-M. A. Comp; Cad / CAM]"
$ws.Range("E14").Value = "[Humberto-Comam. Pneumáticos, Valmir-Calderaria, Gisele-Ens. Dest. não Dest., Cleidson-Metrologia 2]"
$ws.Range("F14").Value = "-"

# Row 15
$ws.Range("B15").Value = "[Ivan-Tec. Soldagem, Wellington-Trat. Térmicos, And. Edson-Usin. CNC, Ludoff-Camam. Hidráulicos]"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "[Joel Lima-Tec. Fundição, Aderci-Fresagem, Anselmo-M. Motor Endot., Ludoff-Cont. Lóg. Prog. CLP]"
$ws.Range("E15").Value = "[Aderci-M. A. Comp; Cad / CAM, -, -, Leonardo R.-Mec. Manut. Equip. Ind.]"
$ws.Range("F15").Value = "-"

# Row 16
$ws.Range("B16").Value = "[-, Wellington-Trat. Térmicos, Leonardo R.-Mec. Manut. Equip. Ind., Ludoff-Camam. Hidráulicos]"
$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "[Joel Lima-Tec. Fundição, Aderci-Fresagem, Ivan-Tec. Soldagem, Wellington-Trat. Térmicos]"
$ws.Range("E16").Value = "[Aderci-M. A. Comp; Cad / CAM, Leonardo R.-Retificação, -, -]"
